# Insert a new "Industry" column between "Stock Name" (B) and "Mutual Fund" (C),
# shifting the existing Mutual Fund/Status/Jan_2026/Dec_2025/Oct_2025/MoM/QoQ
# columns one place to the right (C:I -> D:J), then populate the new Industry
# column with the industry classification for each holding.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C; this shifts existing C:I -> D:J
# and automatically carries over styles/formatting.
$ws.Columns("C").Insert()

# Header
$ws.Range("C1").Value = "Industry"

# Industry values per row (ISIN / Stock Name order matches the existing rows)
$industries = @{
    2  = "Construction"
    3  = "Banks"
    4  = "Metals & Minerals Trading"
    5  = "Power"
    6  = "Finance"
    7  = "Banks"
    8  = "Auto Components"
    9  = "Insurance"
    10 = "Banks"
    11 = "Automobiles"
    12 = "Banks"
    13 = "Insurance"
    14 = "Power"
    15 = "Realty"
    16 = "Pharmaceuticals & Biotechnology"
    17 = "IT - Software"
    18 = "IT - Software"
    19 = "Insurance"
    20 = "Banks"
    21 = "Pharmaceuticals & Biotechnology"
    22 = "Finance"
}

foreach ($row in $industries.Keys) {
    $ws.Cells.Item($row, 3).Value = $industries[$row]
}
